# Catalog.xlsx update: add UNSPSC / Supplier / Manufacturer columns to the
# SmartForm sheet, consolidate Quantity / Unit of Measure / Price into
# multi-value cells, and flip which sheet/tab is active.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # SmartForm
$ws2 = $wb.Worksheets.Item(2)   # VerifyCatalogSearch

# ---------------------------------------------------------------------
# SmartForm: insert two new columns after "Item Description" (D,E) and
# two more after "Sub Category" (which will then sit at columns I,J).
# Using Insert() (rather than overwriting in place) lets the columns
# that don't change meaning - Category Type / Category / Sub Category -
# carry their original widths along for the ride.
# ---------------------------------------------------------------------
$ws1.Columns.Item(4).Insert()
$ws1.Columns.Item(4).Insert()
$ws1.Columns.Item(9).Insert()
$ws1.Columns.Item(9).Insert()

# Header row
$ws1.Range("A1").Value = "Role"
$ws1.Range("B1").Value = "Location"
$ws1.Range("C1").Value = "Item Description "
$ws1.Range("D1").Value = "UNSPSC Code"
$ws1.Range("E1").Value = "Suggested Supplier(s)"
$ws1.Range("F1").Value = "Category Type"
$ws1.Range("G1").Value = "Category"
$ws1.Range("H1").Value = "Sub Category"
$ws1.Range("I1").Value = "Manufacturer Name"
$ws1.Range("J1").Value = "Manufacturer Part Number"
$ws1.Range("K1").Value = "Quantity"
$ws1.Range("L1").Value = "Unit of Measure"
$ws1.Range("M1").Value = "Price "

# Data row
$ws1.Range("A2").Value = "REQUESTOR"
$ws1.Range("B2").Value = "XEEVA -MJ"
$ws1.Range("C2").Value = "REPOFLOR 100 MG"
$ws1.Range("D2").Value = "UNSPSC001"
$ws1.Range("E2").Value = "Sachin Supplier Magna"
$ws1.Range("F2").Value = "INFORMATION TECHNOLOGY"
$ws1.Range("G2").Value = "HARDWARE"
$ws1.Range("H2").Value = "CELL PHONES"
$ws1.Range("I2").Value = "ARMSTRONG"
$ws1.Range("J2").Value = "MPN001"
$ws1.Range("K2").Value = "1;2"
$ws1.Range("L2").Value = "EA-EACH;CU-CUBIC"
$ws1.Range("M2").Value = "1;10"

# Column widths for the brand-new / resized columns (the carried-over
# columns F, G, H keep their original widths automatically).
$ws1.Columns.Item(3).ColumnWidth  = 16.833333   # Item Description
$ws1.Columns.Item(4).ColumnWidth  = 21.166667   # UNSPSC Code
$ws1.Columns.Item(5).ColumnWidth  = 21.166667   # Suggested Supplier(s)
$ws1.Columns.Item(9).ColumnWidth  = 19.333333   # Manufacturer Name
$ws1.Columns.Item(10).ColumnWidth = 19.333333   # Manufacturer Part Number
$ws1.Columns.Item(11).ColumnWidth = 14.5        # Quantity
$ws1.Columns.Item(12).ColumnWidth = 14.666667   # Unit of Measure
$ws1.Columns.Item(13).ColumnWidth = 9.833333    # Price

# ---------------------------------------------------------------------
# VerifyCatalogSearch: values are unchanged, just re-written so the
# sheet stays internally consistent (shared-string table compaction is
# handled automatically once SmartForm no longer references the old
# "Test Item " / "EA-EACH" strings).
# ---------------------------------------------------------------------
$ws2.Range("A1").Value = "Role"
$ws2.Range("B1").Value = "Location"
$ws2.Range("C1").Value = "CatalogSearchInput "
$ws2.Range("D1").Value = "LocalSearch"
$ws2.Range("E1").Value = "GlobalSearch"
$ws2.Range("F1").Value = "BPO"

$ws2.Range("A2").Value = "REQUESTOR"
$ws2.Range("B2").Value = "MAGNA DECOPLAS"
$ws2.Range("C2").Value = "COMPUTER DESKTOP"
$ws2.Range("D2").Value = "local"
$ws2.Range("E2").Value = "global"
$ws2.Range("F2").Value = "bpo"

# ---------------------------------------------------------------------
# View state: SmartForm becomes the selected/active tab (it was
# VerifyCatalogSearch before); update each sheet's remembered selection.
# ---------------------------------------------------------------------
$ws2.Range("G1").Select()
$ws1.Activate()
$ws1.Range("L6").Select()
